$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap A1/B1 back, keep C1:F1, then add new Type columns ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "PokedexNumber"
$ws.Range("C1").Value = "BaseHealth"
$ws.Range("D1").Value = "BaseAttack"
$ws.Range("E1").Value = "BaseDefense"
$ws.Range("F1").Value = "BaseSpeed"

# --- Row 2 (Alien) ---
$ws.Range("A2").Value = "Alien"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 39
$ws.Range("D2").Value = 52
$ws.Range("E2").Value = 43
$ws.Range("F2").Value = 65

# --- Row 3 (Birb) ---
$ws.Range("A3").Value = "Birb"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 40
$ws.Range("D3").Value = 45
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 56

# --- New Type2/Type1 columns, entered H1 before G1 to match authoring order ---
$ws.Range("H1").Value = "Type2"
$ws.Range("G1").Value = "Type1"

# --- New row 4 (Birb-H) name entered before the rest of the new type data ---
$ws.Range("A4").Value = "Birb-H"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 56
$ws.Range("D4").Value = 45
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = 40

# --- Type columns filled row by row ---
$ws.Range("G2").Value = "Fire"
$ws.Range("H2").Value = "None"
$ws.Range("G3").Value = "Normal"
$ws.Range("H3").Value = "Flying"
$ws.Range("G4").Value = "Ice"
$ws.Range("H4").Value = "Flying"

# --- Update selection to match target ---
$ws.Range("G5").Select()
